$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Goal (per diff):
#  1. Insert a new "Meta description" paragraph right after the first
#     (Heading1) paragraph. It reuses the same run layout (leading empty run
#     + bold run + plain run) that the old duplicate heading/description
#     paragraphs near the end of the document already use.
#  2. Remove the duplicate bold paragraph
#     ("Play Deep Sea Magic for Free - Review and Game Features") that used
#     to sit just before the final italic paragraph.
#  3. Replace the text of that final italic paragraph with the new DALL-E
#     image-prompt text, keeping its italic formatting intact.
# ---------------------------------------------------------------------------

# Locate the duplicate bold "title" paragraph near the end of the document
# (paragraph 1, the Heading1, has identical text, so search from the bottom
# up and stop at the very first/last match found).
# (Range.Text includes a trailing paragraph-mark character, so TrimEnd()
# before comparing.)
$totalBefore = $d.Paragraphs.Count
$boldHeadingIndex = -1
for ($i = $totalBefore; $i -ge 2; $i--) {
    if ($boldHeadingIndex -eq -1 -and $d.Paragraphs($i).Range.Text.TrimEnd() -eq "Play Deep Sea Magic for Free - Review and Game Features") {
        $boldHeadingIndex = $i
    }
}

# --- Step 1: create the new paragraph right after paragraph 1 (the H1) ----
$d.Paragraphs(1).Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(2)
$newPara.Style = "Normal"

# Copy the run structure (leading empty run + bold run) from the bold
# duplicate-heading paragraph near the bottom of the document, so the new
# paragraph matches the same "<w:r/><w:r><w:rPr><w:b/></w:rPr>...</w:r>"
# shape, then overwrite its text.
$sourcePara = $d.Paragraphs($boldHeadingIndex + 1)
$copiedFormattedText = $sourcePara.Range.FormattedText
$newPara.Range.FormattedText = $copiedFormattedText

$newPara = $d.Paragraphs(2)
$newParaRange = $newPara.Range
$boldRunRange = $d.Range($newParaRange.Start, $newParaRange.End - 1)
$boldRunRange.Text = "Meta description"

$newPara = $d.Paragraphs(2)
$newParaRange = $newPara.Range
$tailInsertionPoint = $d.Range($newParaRange.End - 1, $newParaRange.End - 1)
$tailInsertionPoint.InsertAfter(": Experience an exciting gameplay with innovative features in Deep Sea Magic slot game. Play for free and read our full review here.")

# --- Step 2: delete the now-shifted duplicate bold heading paragraph ------
$boldHeadingIndex = -1
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    if ($boldHeadingIndex -eq -1 -and $d.Paragraphs($i).Range.Text.TrimEnd() -eq "Play Deep Sea Magic for Free - Review and Game Features") {
        $boldHeadingIndex = $i
    }
}
$d.Paragraphs($boldHeadingIndex).Range.Delete()

# --- Step 3: replace the italic paragraph's text (now the last paragraph) -
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastParaRange = $lastPara.Range
$lastTextRange = $d.Range($lastParaRange.Start, $lastParaRange.End - 1)
$lastTextRange.Text = "DALLE, please create a feature image for `"Deep Sea Magic`" that fits the game's theme and features a happy Maya warrior with glasses in a cartoon style. The image should be eye-catching and playful, incorporating elements of the deep sea and the game's features such as the Drop & Lock feature, bonuses, and wild symbols. Please make sure that the image is high-quality and in line with the overall aesthetic of the game. Thank you!"

Write-Output "Edit complete."
